# Updated symbol list on Wed Dec 28 10:55:23 UTC 2022 with GitHub Actions
#
# The "Price" column (D) stores values as plain text (not numbers), so we
# force each target cell to Text format before writing the new numeric-
# looking string. This prevents Excel's automatic type inference from
# converting the text into a floating point number (which would also
# introduce binary floating-point rounding artifacts, e.g. 243.39 ->
# 243.38999999999999).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
}

Set-TextValue "D2"  "243.39"
Set-TextValue "D4"  "5.265"
Set-TextValue "D5"  "0.05813"
Set-TextValue "D6"  "6.478"
Set-TextValue "D7"  "3.330"
Set-TextValue "D8"  "0.8079"
Set-TextValue "D9"  "0.8738"
Set-TextValue "D10" "0.1386"
Set-TextValue "D13" "0.03052"
Set-TextValue "D14" "0.09305"
Set-TextValue "D15" "3.862"
Set-TextValue "D16" "0.001550"
Set-TextValue "D17" "0.04684"
Set-TextValue "D18" "0.0006030"
$ws.Range("E18").Value = "17OneONE"
Set-TextValue "D19" "0.006168"
Set-TextValue "D20" "0.001261"
Set-TextValue "D21" "0.004592"
Set-TextValue "D22" "0.00008700"
Set-TextValue "D23" "3.560"
Set-TextValue "D24" "2.172"
Set-TextValue "D40" "0.03786"
Set-TextValue "D41" "0.006300"
Set-TextValue "D42" "0.1051"
Set-TextValue "D44" "0.007963"
Set-TextValue "D45" "0.00005544"
Set-TextValue "D47" "0.5660"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINWorstin24h"
Set-TextValue "D48" "0.01402"
